# Generate Report for Handoff
# Adds the new source file "bf64de53-ac12-4f98-b4f7-d6e6011dc4b1.md" (handed off for
# localization) as a new row on every sheet, pushing the ".localization-config" /
# "Not to be localized" row down by one.

$wb = $excel.ActiveWorkbook

$newMd    = "bf64de53-ac12-4f98-b4f7-d6e6011dc4b1.md"
$newZhXlf = "bf64de53-ac12-4f98-b4f7-d6e6011dc4b1.8438a1e037094bca3c3cb27aaec2309b95e6b8d9.zh-cn.xlf"
$newDeXlf = "bf64de53-ac12-4f98-b4f7-d6e6011dc4b1.8438a1e037094bca3c3cb27aaec2309b95e6b8d9.de-de.xlf"
$newZhDt  = "2016-03-10 18:36:00"
$newDeDt  = "2016-03-10 18:36:05"

$newMdUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/c0930bafda350e20443d571f280161b768437c6d/e2e/bf64de53-ac12-4f98-b4f7-d6e6011dc4b1.md"
$cfgUrl      = "https://github.com/OpenLocalizationTest/oltest/blob/ef0e3c4afd08a185d564564b1bb70e25f843ab25/.localization-config"
$newZhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0cab619e0fff4d4bcaf590b1bfab8fd612acb9ed/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newZhXlf"
$newDeXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8bbebc01928feaf0b4b60fad591eee00c9518ad4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newDeXlf"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Push the ".localization-config" row from row 3 down to row 4.
$ov.Range("A4").Value = $ov.Range("A3").Value2
$ov.Range("B4").Value = $ov.Range("B3").Value2
$ov.Range("C4").Value = $ov.Range("C3").Value2

# Write the new row 3 for the handed-off file.
$ov.Range("A3").Value = $newMd
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"

# Rebuild the hyperlinks in order (A2, A3, A4) so relationship ids line up.
$ov.Range("A2").Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c0930bafda350e20443d571f280161b768437c6d/e2e/00cbf9ee-0a9a-4a30-84ad-1ec060825999.md", "", "", "00cbf9ee-0a9a-4a30-84ad-1ec060825999.md")
$ov.Hyperlinks.Add($ov.Range("A3"), $newMdUrl, "", "", $newMd)
$ov.Hyperlinks.Add($ov.Range("A4"), $cfgUrl, "", "", ".localization-config")
$ov.Range("A2:A4").Font.Underline = 2
$ov.Range("A2:A4").Font.Color = 15570276

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Push the ".localization-config" row from row 3 down to row 4.
$zh.Range("A4").Value = $zh.Range("A3").Value2
$zh.Range("B4").Value = $zh.Range("B3").Value2
$zh.Range("D4").Value = $zh.Range("D3").Value2
$zh.Range("G4").Value = $zh.Range("G3").Value2
$zh.Range("H4").Value = "Ignored"

# Write the new row 3 for the handed-off file.
$zh.Range("A3").Value = $newMd
$zh.Range("B3").Value = "Ready for handoff"
$zh.Range("C3").Value = $newZhXlf
$zh.Range("D3").Value = $newZhDt
$zh.Range("G3").Value = "0001-01-01 00:00:00"
$zh.Range("H3").Value = "Include"

$zh.Range("D2:D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Rebuild the hyperlinks in order (A2, C2, A3, C3, A4) so relationship ids line up.
$zh.Range("A2").Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c0930bafda350e20443d571f280161b768437c6d/e2e/00cbf9ee-0a9a-4a30-84ad-1ec060825999.md", "", "", "00cbf9ee-0a9a-4a30-84ad-1ec060825999.md")
$zh.Hyperlinks.Add($zh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0cab619e0fff4d4bcaf590b1bfab8fd612acb9ed/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/00cbf9ee-0a9a-4a30-84ad-1ec060825999.a0e2c0d7e0fed157e66419b3476280ec3bdc3c13.zh-cn.xlf", "", "", "00cbf9ee-0a9a-4a30-84ad-1ec060825999.a0e2c0d7e0fed157e66419b3476280ec3bdc3c13.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A3"), $newMdUrl, "", "", $newMd)
$zh.Hyperlinks.Add($zh.Range("C3"), $newZhXlfUrl, "", "", $newZhXlf)
$zh.Hyperlinks.Add($zh.Range("A4"), $cfgUrl, "", "", ".localization-config")
$zh.Range("A2:A4").Font.Underline = 2
$zh.Range("A2:A4").Font.Color = 15570276
$zh.Range("C2:C3").Font.Underline = 2
$zh.Range("C2:C3").Font.Color = 15570276

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

# Push the ".localization-config" row from row 3 down to row 4.
$de.Range("A4").Value = $de.Range("A3").Value2
$de.Range("B4").Value = $de.Range("B3").Value2
$de.Range("D4").Value = $de.Range("D3").Value2
$de.Range("G4").Value = $de.Range("G3").Value2
$de.Range("H4").Value = "Ignored"

# Write the new row 3 for the handed-off file.
$de.Range("A3").Value = $newMd
$de.Range("B3").Value = "Ready for handoff"
$de.Range("C3").Value = $newDeXlf
$de.Range("D3").Value = $newDeDt
$de.Range("G3").Value = "0001-01-01 00:00:00"
$de.Range("H3").Value = "Include"

$de.Range("D2:D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Rebuild the hyperlinks in order (A2, C2, A3, C3, A4) so relationship ids line up.
$de.Range("A2").Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c0930bafda350e20443d571f280161b768437c6d/e2e/00cbf9ee-0a9a-4a30-84ad-1ec060825999.md", "", "", "00cbf9ee-0a9a-4a30-84ad-1ec060825999.md")
$de.Hyperlinks.Add($de.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8bbebc01928feaf0b4b60fad591eee00c9518ad4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/00cbf9ee-0a9a-4a30-84ad-1ec060825999.a0e2c0d7e0fed157e66419b3476280ec3bdc3c13.de-de.xlf", "", "", "00cbf9ee-0a9a-4a30-84ad-1ec060825999.a0e2c0d7e0fed157e66419b3476280ec3bdc3c13.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A3"), $newMdUrl, "", "", $newMd)
$de.Hyperlinks.Add($de.Range("C3"), $newDeXlfUrl, "", "", $newDeXlf)
$de.Hyperlinks.Add($de.Range("A4"), $cfgUrl, "", "", ".localization-config")
$de.Range("A2:A4").Font.Underline = 2
$de.Range("A2:A4").Font.Color = 15570276
$de.Range("C2:C3").Font.Underline = 2
$de.Range("C2:C3").Font.Color = 15570276

Write-Output "Report for Handoff generated"
